$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column B (shifts existing B:R to C:S)
$ws.Columns("B:B").Insert()

# Populate the new "Usertype" column
$ws.Range("B1").Value = "Usertype"
$ws.Range("B2").Value = "Admin User"
$ws.Range("B7").Value = "Staff User"
$ws.Range("B12").Value = "Client User"

# Give the new column the same width as column A (no auto-fit applied to it)
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Update the view: move selection to where the user left off
$ws.Range("C12").Select()
